$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-27 with the new weekly data values ---
# Row 2
$ws.Range("D2").Value = 44518
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14571
$ws.Range("O2").Value = 'Región del Maule'
$ws.Range("P2").Value = 583

# Row 3
$ws.Range("D3").Value = 44629
$ws.Range("J3").Value = 35
$ws.Range("K3").Value = 25000
$ws.Range("L3").Value = 26000
$ws.Range("M3").Value = 25429
$ws.Range("N3").Value = '$/saco 25 kilos'
$ws.Range("O3").Value = 'Región Metropolitana'
$ws.Range("P3").Value = 1017

# Row 4
$ws.Range("D4").Value = 44496
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14520
$ws.Range("N4").Value = '$/malla 25 kilos'
$ws.Range("O4").Value = 'Provincia de Huasco'
$ws.Range("P4").Value = 581

# Row 5
$ws.Range("D5").Value = 44482
$ws.Range("H5").Value = 'Perfection'
$ws.Range("J5").Value = 130
$ws.Range("K5").Value = 24000
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = 24385
$ws.Range("O5").Value = 'Región de O''Higgins'
$ws.Range("P5").Value = 975

# Row 6
$ws.Range("D6").Value = 44505
$ws.Range("J6").Value = 210
$ws.Range("K6").Value = 6500
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 6714
$ws.Range("O6").Value = 'Región del Maule'
$ws.Range("P6").Value = 269

# Row 7
$ws.Range("D7").Value = 44657
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 24000
$ws.Range("L7").Value = 25000
$ws.Range("M7").Value = 24400
$ws.Range("P7").Value = 976

# Row 8
$ws.Range("D8").Value = 44631
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 24000
$ws.Range("L8").Value = 25000
$ws.Range("M8").Value = 24467
$ws.Range("O8").Value = 'Carahue'
$ws.Range("P8").Value = 979

# Row 9
$ws.Range("D9").Value = 44503
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 15500
$ws.Range("N9").Value = '$/malla 25 kilos'
$ws.Range("O9").Value = 'Provincia de Limarí'
$ws.Range("P9").Value = 620

# Row 10
$ws.Range("D10").Value = 44539
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 14000
$ws.Range("M10").Value = 13400
$ws.Range("P10").Value = 536

# Row 11
$ws.Range("D11").Value = 44545
$ws.Range("J11").Value = 180
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 16000
$ws.Range("M11").Value = 15444
$ws.Range("O11").Value = 'Carahue'
$ws.Range("P11").Value = 618

# Row 13
$ws.Range("D13").Value = 44589
$ws.Range("J13").Value = 160
$ws.Range("K13").Value = 22000
$ws.Range("L13").Value = 23000
$ws.Range("M13").Value = 22500
$ws.Range("N13").Value = '$/malla 25 kilos'
$ws.Range("P13").Value = 900

# Row 14
$ws.Range("D14").Value = 44454
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 36000
$ws.Range("L14").Value = 38000
$ws.Range("M14").Value = 37000
$ws.Range("P14").Value = 1480

# Row 15
$ws.Range("D15").Value = 44519
$ws.Range("J15").Value = 240
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 17583
$ws.Range("P15").Value = 703

# Row 16
$ws.Range("D16").Value = 44328
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 33000
$ws.Range("L16").Value = 34000
$ws.Range("M16").Value = 33500
$ws.Range("N16").Value = '$/malla 25 kilos'
$ws.Range("O16").Value = 'Provincia de Huasco'
$ws.Range("P16").Value = 1340

# Row 17
$ws.Range("D17").Value = 44532
$ws.Range("H17").Value = 'Sin especificar'
$ws.Range("J17").Value = 250
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14400
$ws.Range("N17").Value = '$/saco 25 kilos'
$ws.Range("O17").Value = 'Región del Maule'
$ws.Range("P17").Value = 576

# Row 18
$ws.Range("D18").Value = 44399
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 39000
$ws.Range("L18").Value = 40000
$ws.Range("M18").Value = 39600
$ws.Range("P18").Value = 1584

# Row 19
$ws.Range("D19").Value = 44615
$ws.Range("H19").Value = 'Sin especificar'
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = 28000
$ws.Range("L19").Value = 30000
$ws.Range("M19").Value = 29000
$ws.Range("N19").Value = '$/saco 25 kilos'
$ws.Range("O19").Value = 'Carahue'
$ws.Range("P19").Value = 1160

# Row 20
$ws.Range("D20").Value = 44643
$ws.Range("J20").Value = 90
$ws.Range("K20").Value = 25000
$ws.Range("L20").Value = 26000
$ws.Range("M20").Value = 25444
$ws.Range("O20").Value = 'Carahue'
$ws.Range("P20").Value = 1018

# Row 21
$ws.Range("D21").Value = 44483
$ws.Range("J21").Value = 220
$ws.Range("K21").Value = 19000
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = 19455
$ws.Range("N21").Value = '$/saco 25 kilos'
$ws.Range("O21").Value = 'Región Metropolitana'
$ws.Range("P21").Value = 778

# Row 22
$ws.Range("D22").Value = 44533
$ws.Range("H22").Value = 'Perfection'
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 14000
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = 14375
$ws.Range("N22").Value = '$/malla 25 kilos'
$ws.Range("P22").Value = 575

# Row 23
$ws.Range("D23").Value = 44595
$ws.Range("H23").Value = 'Perfection'
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 26000
$ws.Range("L23").Value = 28000
$ws.Range("M23").Value = 27200
$ws.Range("P23").Value = 1088

# Row 24
$ws.Range("D24").Value = 44342
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = 30000
$ws.Range("L24").Value = 32000
$ws.Range("M24").Value = 31000
$ws.Range("O24").Value = 'Provincia de Limarí'
$ws.Range("P24").Value = 1240

# Row 25
$ws.Range("D25").Value = 44162
$ws.Range("H25").Value = 'Sin especificar'
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 17000
$ws.Range("L25").Value = 18000
$ws.Range("M25").Value = 17500
$ws.Range("O25").Value = 'Región del Maule'
$ws.Range("P25").Value = 700

# Row 26
$ws.Range("D26").Value = 44512
$ws.Range("J26").Value = 100
$ws.Range("M26").Value = 14500
$ws.Range("N26").Value = '$/saco 25 kilos'
$ws.Range("P26").Value = 580

# Row 27
$ws.Range("D27").Value = 44517
$ws.Range("J27").Value = 110
$ws.Range("K27").Value = 17000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 17455
$ws.Range("N27").Value = '$/saco 25 kilos'
$ws.Range("O27").Value = 'Región del Maule'
$ws.Range("P27").Value = 698

# --- Append new row 28 ---
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = 'Vega Monumental Concepción'
$ws.Range("C28").Value = 'Bíobío'
$ws.Range("D28").Value = 44335
$ws.Range("E28").Value = 8
$ws.Range("F28").Value = 100112022
$ws.Range("G28").Value = 'Arveja Verde'
$ws.Range("H28").Value = 'Perfection'
$ws.Range("I28").Value = 'Primera'
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 30000
$ws.Range("L28").Value = 32000
$ws.Range("M28").Value = 31000
$ws.Range("N28").Value = '$/malla 25 kilos'
$ws.Range("O28").Value = 'Provincia de Huasco'
$ws.Range("P28").Value = 1240
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = 'Hortaliza'

# Preserve the datetime number format on the new date cell, matching column D style
$ws.Range("D28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
